$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Murtaza's password hash (E2)
$ws.Range("E2").Value = '$2b$10$S.yjhL0mAORMSudaq7QcmupyZfO7xD13Oxdxz1ZHqxLAijLVDeaQ6'

# Append a new record for Naeem Shaikh on row 4
$ws.Range("A4").Value = 7
$ws.Range("B4").Value = "Naeem"
$ws.Range("C4").Value = "Shaikh"
$ws.Range("D4").Value = "chaseyourdreams121@gmail.com"
$ws.Range("E4").Value = '$2b$10$jQHNu2bsZwmZWz5kQBYmce33Ps8KvDIWwCdgdpWJkXNnwIzIeIoP.'
